$d = $word.ActiveDocument

# Paragraph 1: "Các qui định về họp hành nội bộ" -> append italic "." as a new run
$rng1 = $d.Content
$rng1.Find.Execute("ọp hành nội bộ") | Out-Null
$rng1.Collapse(0)
$rng1.InsertAfter(".")
$rng1.Font.Italic = $true
$rng1.Font.ItalicBi = $true

# Paragraph 2: "Các qui định về họp hành với khách hàng" -> append italic "." as a new run
$rng2 = $d.Content
$rng2.Find.Execute("ọp hành với khách hàng") | Out-Null
$rng2.Collapse(0)
$rng2.InsertAfter(".")
$rng2.Font.Italic = $true
$rng2.Font.ItalicBi = $true
